# This script applies the "updated names of sigma index" edit to the
# meta_data_names workbook:
#   - F15 "Estimate_sigma_index" -> "Estimate_sigma_survey"
#   - F16 "Sigma_index_prior"    -> "Sigma_survey_prior"
#   - The three NOTE rows that used to live in F25:F27 are moved down to
#     F29:F31 (F25:F27 are cleared), and the last NOTE's text gains a
#     trailing clause.
#   - The sheet view's scroll position / active selection is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sigma index labels (column F) ---------------------------
$ws.Range("F15").Value = "Estimate_sigma_survey"
$ws.Range("F16").Value = "Sigma_survey_prior"

# --- Move the trailing NOTE column from rows 25-27 down to rows 29-31 ---
$ws.Range("F25:F27").Cut($ws.Range("F29:F31"))

# Update the moved note's wording
$ws.Range("F31").Value = "NOTE: Columns for ages are index by 1 trhough nages, but are place holders."

# --- Update sheet view scroll/selection state ----------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C29").Select()
